$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: change the monthly values from 17500 to 4560 (B4:M4). N4 keeps its SUM formula.
$ws.Range("B4:M4").Value = 4560

# Row 5: replace literal values with SUM formulas referencing rows 2:4 of the same column.
# B5 gets its own (non-shared) formula; C5:N5 is filled as one shared-formula group.
$ws.Range("B5").Formula = "=SUM(B2:B4)"
$ws.Range("C5:N5").Formula = "=SUM(C2:C4)"

# Row 14: re-fill the existing totals formula across the row so it becomes a shared
# formula group (same formula/results as before, just stored more compactly).
$ws.Range("B14:M14").Formula = "=SUM(B12:B16)"

# Update the selected cell shown in the sheet view.
$ws.Range("L10").Select()
